$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.925.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.052.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.22'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.51%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.386'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0809'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.356.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.042.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.807.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0830'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E28').Value = '  +5.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.54'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.06'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.02%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0612'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.39%  '
$ws.Range('E38').Value = '  +4.92%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.486.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0935'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +17.42%  '
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.244.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.44%  '
